$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to reflect new "damage" naming scheme
$ws.Range("F1").Value = "Body.damage"
$ws.Range("G1").Value = "Forewing.dorsal.damage"
$ws.Range("H1").Value = "Forewing.ventral.damage"
$ws.Range("I1").Value = "Hindwing.dorsal.damage"
$ws.Range("J1").Value = "Hindwing.ventral.damage"

# Fix typo: "Speced" -> "Spectra"
$ws.Range("R1").Value = "Spectra"
